$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 235, shifting existing rows 235..349 down to 236..350
$ws.Rows(235).Insert()

# Populate the newly inserted row 235 with the new record's data
$ws.Cells.Item(235, 1).Value = 6
$ws.Cells.Item(235, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(235, 3).Value = 'Metropolitana'
$ws.Cells.Item(235, 4).Value = 44466
$ws.Cells.Item(235, 5).Value = 13
$ws.Cells.Item(235, 6).Value = 100112044
$ws.Cells.Item(235, 7).Value = 'Perejil'
$ws.Cells.Item(235, 8).Value = 'Sin especificar'
$ws.Cells.Item(235, 9).Value = 'Primera'
$ws.Cells.Item(235, 10).Value = 140
$ws.Cells.Item(235, 11).Value = 7500
$ws.Cells.Item(235, 12).Value = 8000
$ws.Cells.Item(235, 13).Value = 7714
$ws.Cells.Item(235, 14).Value = '$/docena de atados'
$ws.Cells.Item(235, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(235, 16).Value = 2571
$ws.Cells.Item(235, 17).Value = 3
$ws.Cells.Item(235, 18).Value = 'Hortaliza'
